$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "flash mode" (row 7 - take photo flash option) and "app icon" related
# items as completed ("y") in column C, matching the formatting already used
# by the existing "y" cells (e.g. C3).
$cells = @("C2", "C4", "C5", "C7")
foreach ($addr in $cells) {
    $cell = $ws.Range($addr)
    $cell.Value = "y"
    $cell.HorizontalAlignment = -4108
    $cell.Font.Size = 18
}

# Move the active selection.
$ws.Range("G10").Select()
